$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "time_taken" in F1, copying the header style/format from E1 (s="1")
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F1").Value = "time_taken"

# Fill F2:F105 with the recorded time_taken values for each data row
$timeTaken = New-Object 'object[,]' 104,1
$timeTaken[0,0] = '2021-10-05 13:40:15.069030'
$timeTaken[1,0] = '2021-10-05 13:40:15.069043'
$timeTaken[2,0] = '2021-10-05 13:40:15.069047'
$timeTaken[3,0] = '2021-10-05 13:40:15.069051'
$timeTaken[4,0] = '2021-10-05 13:40:15.069054'
$timeTaken[5,0] = '2021-10-05 13:40:15.069057'
$timeTaken[6,0] = '2021-10-05 13:40:15.069060'
$timeTaken[7,0] = '2021-10-05 13:40:15.069063'
$timeTaken[8,0] = '2021-10-05 13:40:15.069066'
$timeTaken[9,0] = '2021-10-05 13:40:15.069069'
$timeTaken[10,0] = '2021-10-05 13:40:15.069072'
$timeTaken[11,0] = '2021-10-05 13:40:15.069075'
$timeTaken[12,0] = '2021-10-05 13:40:15.069079'
$timeTaken[13,0] = '2021-10-05 13:40:15.069081'
$timeTaken[14,0] = '2021-10-05 13:40:15.069085'
$timeTaken[15,0] = '2021-10-05 13:40:15.069087'
$timeTaken[16,0] = '2021-10-05 13:40:15.069091'
$timeTaken[17,0] = '2021-10-05 13:40:15.069094'
$timeTaken[18,0] = '2021-10-05 13:40:15.069097'
$timeTaken[19,0] = '2021-10-05 13:40:15.069100'
$timeTaken[20,0] = '2021-10-05 13:40:15.069103'
$timeTaken[21,0] = '2021-10-05 13:40:15.069106'
$timeTaken[22,0] = '2021-10-05 13:40:15.069109'
$timeTaken[23,0] = '2021-10-05 13:40:15.069112'
$timeTaken[24,0] = '2021-10-05 13:40:15.069115'
$timeTaken[25,0] = '2021-10-05 13:40:15.069119'
$timeTaken[26,0] = '2021-10-05 13:40:15.069121'
$timeTaken[27,0] = '2021-10-05 13:40:15.069124'
$timeTaken[28,0] = '2021-10-05 13:40:15.069127'
$timeTaken[29,0] = '2021-10-05 13:40:15.069130'
$timeTaken[30,0] = '2021-10-05 13:40:15.069133'
$timeTaken[31,0] = '2021-10-05 13:40:15.069136'
$timeTaken[32,0] = '2021-10-05 13:40:15.069140'
$timeTaken[33,0] = '2021-10-05 13:40:15.069143'
$timeTaken[34,0] = '2021-10-05 13:40:15.069146'
$timeTaken[35,0] = '2021-10-05 13:40:15.069149'
$timeTaken[36,0] = '2021-10-05 13:40:15.069152'
$timeTaken[37,0] = '2021-10-05 13:40:15.069155'
$timeTaken[38,0] = '2021-10-05 13:40:15.069158'
$timeTaken[39,0] = '2021-10-05 13:40:15.069161'
$timeTaken[40,0] = '2021-10-05 13:40:15.069165'
$timeTaken[41,0] = '2021-10-05 13:40:15.069168'
$timeTaken[42,0] = '2021-10-05 13:40:15.069171'
$timeTaken[43,0] = '2021-10-05 13:40:15.069174'
$timeTaken[44,0] = '2021-10-05 13:40:15.069177'
$timeTaken[45,0] = '2021-10-05 13:40:15.069180'
$timeTaken[46,0] = '2021-10-05 13:40:15.069183'
$timeTaken[47,0] = '2021-10-05 13:40:15.069186'
$timeTaken[48,0] = '2021-10-05 13:40:15.069189'
$timeTaken[49,0] = '2021-10-05 13:40:15.069192'
$timeTaken[50,0] = '2021-10-05 13:40:15.069195'
$timeTaken[51,0] = '2021-10-05 13:40:15.069198'
$timeTaken[52,0] = '2021-10-05 13:40:15.069202'
$timeTaken[53,0] = '2021-10-05 13:40:15.069205'
$timeTaken[54,0] = '2021-10-05 13:40:15.069208'
$timeTaken[55,0] = '2021-10-05 13:40:15.069211'
$timeTaken[56,0] = '2021-10-05 13:40:15.069214'
$timeTaken[57,0] = '2021-10-05 13:40:15.069217'
$timeTaken[58,0] = '2021-10-05 13:40:15.069220'
$timeTaken[59,0] = '2021-10-05 13:40:15.069223'
$timeTaken[60,0] = '2021-10-05 13:40:15.069226'
$timeTaken[61,0] = '2021-10-05 13:40:15.069229'
$timeTaken[62,0] = '2021-10-05 13:40:15.069232'
$timeTaken[63,0] = '2021-10-05 13:40:15.069235'
$timeTaken[64,0] = '2021-10-05 13:40:15.069239'
$timeTaken[65,0] = '2021-10-05 13:40:15.069243'
$timeTaken[66,0] = '2021-10-05 13:40:15.069246'
$timeTaken[67,0] = '2021-10-05 13:40:15.069249'
$timeTaken[68,0] = '2021-10-05 13:40:15.069252'
$timeTaken[69,0] = '2021-10-05 13:40:15.069255'
$timeTaken[70,0] = '2021-10-05 13:40:15.069258'
$timeTaken[71,0] = '2021-10-05 13:40:15.069261'
$timeTaken[72,0] = '2021-10-05 13:40:15.069264'
$timeTaken[73,0] = '2021-10-05 13:40:15.069267'
$timeTaken[74,0] = '2021-10-05 13:40:15.069270'
$timeTaken[75,0] = '2021-10-05 13:40:15.069273'
$timeTaken[76,0] = '2021-10-05 13:40:15.069279'
$timeTaken[77,0] = '2021-10-05 13:40:15.069282'
$timeTaken[78,0] = '2021-10-05 13:40:15.069285'
$timeTaken[79,0] = '2021-10-05 13:40:15.069289'
$timeTaken[80,0] = '2021-10-05 13:40:15.069292'
$timeTaken[81,0] = '2021-10-05 13:40:15.069295'
$timeTaken[82,0] = '2021-10-05 13:40:15.069298'
$timeTaken[83,0] = '2021-10-05 13:40:15.069301'
$timeTaken[84,0] = '2021-10-05 13:40:15.069304'
$timeTaken[85,0] = '2021-10-05 13:40:15.069307'
$timeTaken[86,0] = '2021-10-05 13:40:15.069310'
$timeTaken[87,0] = '2021-10-05 13:40:15.069313'
$timeTaken[88,0] = '2021-10-05 13:40:15.069316'
$timeTaken[89,0] = '2021-10-05 13:40:15.069319'
$timeTaken[90,0] = '2021-10-05 13:40:15.069322'
$timeTaken[91,0] = '2021-10-05 13:40:15.069325'
$timeTaken[92,0] = '2021-10-05 13:40:15.069330'
$timeTaken[93,0] = '2021-10-05 13:40:15.069333'
$timeTaken[94,0] = '2021-10-05 13:40:15.069337'
$timeTaken[95,0] = '2021-10-05 13:40:15.069340'
$timeTaken[96,0] = '2021-10-05 13:40:15.069343'
$timeTaken[97,0] = '2021-10-05 13:40:15.069346'
$timeTaken[98,0] = '2021-10-05 13:40:15.069349'
$timeTaken[99,0] = '2021-10-05 13:40:15.069352'
$timeTaken[100,0] = '2021-10-05 13:40:15.069355'
$timeTaken[101,0] = '2021-10-05 13:40:15.069358'
$timeTaken[102,0] = '2021-10-05 13:40:15.069361'
$timeTaken[103,0] = '2021-10-05 13:40:15.069364'
$ws.Range("F2:F105").Value = $timeTaken

Write-Output "time_taken column added"
